$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column B ("PDL Date") values from "01012025" to "99999999"
# for all data rows (2 through 383), leaving the header row (B1) untouched.
# Force the cells to remain text (not get auto-converted to numbers),
# then restore the default "Normal" style so no stray formatting change
# is introduced on the cells.
$lastRow = 383
$range = $ws.Range("B2:B" + $lastRow)
$range.NumberFormat = "@"
$range.Value = "99999999"
$range.Style = "Normal"
